# Update "想去人数" (interest count) figures in column F across the
# workbook's sheets to the freshly scraped values (gh-pages data refresh,
# commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$updates = @{
    2  = 8406
    3  = 36702
    5  = 618
    9  = 455
    10 = 825
    11 = 81
    12 = 664
    13 = 500
    15 = 610
    17 = 452
    18 = 436
    19 = 1140
    21 = 782
    22 = 2446
    23 = 941
    24 = 536
    26 = 1131
    28 = 718
    29 = 718
}
foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$updates = @{
    4  = 361
    7  = 54
    9  = 140
    12 = 8
}
foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$updates = @{
    2 = 595
}
foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$updates = @{
    2  = 595
    3  = 8406
    5  = 36702
    7  = 618
    12 = 455
    13 = 361
    16 = 825
    17 = 81
    18 = 664
    19 = 500
    20 = 54
    23 = 140
    26 = 610
    28 = 452
    29 = 436
    30 = 1140
    32 = 782
    33 = 2446
    34 = 941
    35 = 536
    37 = 1131
    39 = 8
    40 = 718
    41 = 718
}
foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
